# ChatTrinityDramaSourceData.xlsx update
# - Insert a new "Web" data-source row at row 47 (pushing the existing
#   "PDF" rows down by one, 47-76 -> 48-77) pointing at the Trinity
#   "registered exam centre" page, formatted as a real hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new row, shifting everything currently at/after
# row 47 down by one.
$ws.Rows.Item(47).Insert()

# Column A keeps using the existing "Web" label already used by rows 2-46.
$ws.Range("A47").Value = "Web"

# Column B becomes a clickable hyperlink to the new source; Excel applies
# the built-in "Hyperlink" cell style (blue + underline) automatically.
$ws.Hyperlinks.Add($ws.Range("B47"), "https://www.trinitycollege.com/about-us/work-with-trinity/registered-exam-centre")
